$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H26").Value = 7.36
$ws.Range("I26").Value = 86.51000000000001
$ws.Range("J26").Value = 155.12
$ws.Range("K26").Value = 195.12
$ws.Range("M26").Value = 587.14
$ws.Range("H27").Value = 7.36
$ws.Range("I27").Value = 86.51000000000001
$ws.Range("J27").Value = 155.12
$ws.Range("K27").Value = 195.12
$ws.Range("M27").Value = 587.14
$ws.Range("H28").Value = 7.36
$ws.Range("I28").Value = 86.38
$ws.Range("J28").Value = 155.1
$ws.Range("K28").Value = 195.09
$ws.Range("M28").Value = 587.03
$ws.Range("H29").Value = 19.95
$ws.Range("I29").Value = 84.43000000000001
$ws.Range("J29").Value = 143.97
$ws.Range("K29").Value = 193.19
$ws.Range("M29").Value = 614.8
$ws.Range("H30").Value = 7.36
$ws.Range("I30").Value = 86.38
$ws.Range("J30").Value = 155.1
$ws.Range("K30").Value = 195.09
$ws.Range("M30").Value = 587.03
$ws.Range("H31").Value = 19.95
$ws.Range("I31").Value = 84.43000000000001
$ws.Range("J31").Value = 143.97
$ws.Range("K31").Value = 193.18
$ws.Range("M31").Value = 614.8
$ws.Range("H32").Value = 19.98
$ws.Range("I32").Value = 84.01000000000001
$ws.Range("J32").Value = 143.97
$ws.Range("K32").Value = 193.14
$ws.Range("M32").Value = 614.08
$ws.Range("H33").Value = 19.99
$ws.Range("I33").Value = 84.01000000000001
$ws.Range("J33").Value = 143.97
$ws.Range("K33").Value = 193.14
$ws.Range("M33").Value = 614.0599999999999
$ws.Range("H34").Value = 5.81
$ws.Range("I34").Value = 68.34999999999999
$ws.Range("J34").Value = 157.58
$ws.Range("K34").Value = 197.52
$ws.Range("M34").Value = 505.62
$ws.Range("H35").Value = 5.81
$ws.Range("I35").Value = 68.34999999999999
$ws.Range("J35").Value = 157.58
$ws.Range("K35").Value = 197.52
$ws.Range("M35").Value = 505.62
$ws.Range("H36").Value = 5.81
$ws.Range("I36").Value = 68.2
$ws.Range("J36").Value = 157.58
$ws.Range("K36").Value = 197.52
$ws.Range("M36").Value = 505.62
$ws.Range("H37").Value = 12.55
$ws.Range("I37").Value = 53.1
$ws.Range("J37").Value = 157.58
$ws.Range("K37").Value = 204.87
$ws.Range("M37").Value = 577.77
$ws.Range("H38").Value = 5.81
$ws.Range("I38").Value = 68.2
$ws.Range("J38").Value = 157.58
$ws.Range("K38").Value = 197.52
$ws.Range("M38").Value = 505.62
$ws.Range("H39").Value = 12.55
$ws.Range("I39").Value = 53.1
$ws.Range("J39").Value = 157.58
$ws.Range("K39").Value = 204.87
$ws.Range("M39").Value = 577.77
$ws.Range("H40").Value = 12.55
$ws.Range("I40").Value = 52.74
$ws.Range("J40").Value = 157.58
$ws.Range("K40").Value = 204.87
$ws.Range("M40").Value = 577.77
$ws.Range("H41").Value = 12.55
$ws.Range("I41").Value = 52.73
$ws.Range("J41").Value = 157.58
$ws.Range("K41").Value = 204.87
$ws.Range("M41").Value = 577.77
$ws.Range("J42").Value = 216.63
$ws.Range("K42").Value = 234.58
$ws.Range("M42").Value = 597.8
$ws.Range("J43").Value = 216.63
$ws.Range("K43").Value = 234.58
$ws.Range("M43").Value = 597.8
$ws.Range("J44").Value = 216.1
$ws.Range("K44").Value = 234.33
$ws.Range("M44").Value = 597.71
$ws.Range("J45").Value = 162.75
$ws.Range("K45").Value = 218.01
$ws.Range("M45").Value = 605.87
$ws.Range("J46").Value = 216.1
$ws.Range("K46").Value = 234.33
$ws.Range("M46").Value = 597.71
$ws.Range("J47").Value = 162.57
$ws.Range("K47").Value = 217.99
$ws.Range("M47").Value = 605.87
$ws.Range("J48").Value = 161.21
$ws.Range("K48").Value = 217.47
$ws.Range("M48").Value = 605.6
$ws.Range("J49").Value = 161.07
$ws.Range("K49").Value = 217.45
$ws.Range("M49").Value = 605.5700000000001
$ws.Range("H74").Value = 7.72
$ws.Range("I74").Value = 90.7
$ws.Range("J74").Value = 173.8
$ws.Range("K74").Value = 224.08
$ws.Range("M74").Value = 657.92
$ws.Range("H75").Value = 7.72
$ws.Range("I75").Value = 90.7
$ws.Range("J75").Value = 173.8
$ws.Range("K75").Value = 224.08
$ws.Range("M75").Value = 657.92
$ws.Range("H76").Value = 7.72
$ws.Range("I76").Value = 90.56
$ws.Range("J76").Value = 173.82
$ws.Range("K76").Value = 224.09
$ws.Range("M76").Value = 657.4
$ws.Range("H77").Value = 21.34
$ws.Range("I77").Value = 90.34999999999999
$ws.Range("J77").Value = 182.21
$ws.Range("K77").Value = 260.68
$ws.Range("M77").Value = 749.5599999999999
$ws.Range("H78").Value = 7.72
$ws.Range("I78").Value = 90.56
$ws.Range("J78").Value = 173.82
$ws.Range("K78").Value = 224.09
$ws.Range("M78").Value = 657.4
$ws.Range("H79").Value = 21.35
$ws.Range("I79").Value = 90.34999999999999
$ws.Range("J79").Value = 182.21
$ws.Range("K79").Value = 260.67
$ws.Range("M79").Value = 749.54
$ws.Range("H80").Value = 21.38
$ws.Range("I80").Value = 89.88
$ws.Range("J80").Value = 182.3
$ws.Range("K80").Value = 260.67
$ws.Range("M80").Value = 749.39
$ws.Range("H81").Value = 21.38
$ws.Range("I81").Value = 89.88
$ws.Range("J81").Value = 182.29
$ws.Range("K81").Value = 260.67
$ws.Range("M81").Value = 749.38
$ws.Range("H82").Value = 6.08
$ws.Range("I82").Value = 71.44
$ws.Range("J82").Value = 177.55
$ws.Range("K82").Value = 222.66
$ws.Range("M82").Value = 545.48
$ws.Range("H83").Value = 6.08
$ws.Range("I83").Value = 71.44
$ws.Range("J83").Value = 177.55
$ws.Range("K83").Value = 222.66
$ws.Range("M83").Value = 545.48
$ws.Range("H84").Value = 6.08
$ws.Range("I84").Value = 71.29000000000001
$ws.Range("J84").Value = 177.55
$ws.Range("K84").Value = 222.66
$ws.Range("M84").Value = 545.48
$ws.Range("H85").Value = 13.17
$ws.Range("I85").Value = 55.76
$ws.Range("J85").Value = 185.61
$ws.Range("K85").Value = 249.19
$ws.Range("M85").Value = 747.66
$ws.Range("H86").Value = 6.08
$ws.Range("I86").Value = 71.29000000000001
$ws.Range("J86").Value = 177.55
$ws.Range("K86").Value = 222.66
$ws.Range("M86").Value = 545.48
$ws.Range("H87").Value = 13.17
$ws.Range("I87").Value = 55.76
$ws.Range("J87").Value = 185.61
$ws.Range("K87").Value = 249.19
$ws.Range("M87").Value = 747.66
$ws.Range("H88").Value = 13.17
$ws.Range("I88").Value = 55.38
$ws.Range("J88").Value = 185.61
$ws.Range("K88").Value = 249.19
$ws.Range("M88").Value = 747.66
$ws.Range("H89").Value = 13.17
$ws.Range("I89").Value = 55.38
$ws.Range("J89").Value = 185.61
$ws.Range("K89").Value = 249.19
$ws.Range("M89").Value = 747.66
$ws.Range("J90").Value = 229.68
$ws.Range("K90").Value = 245.63
$ws.Range("M90").Value = 607.5700000000001
$ws.Range("J91").Value = 229.68
$ws.Range("K91").Value = 245.63
$ws.Range("M91").Value = 607.5700000000001
$ws.Range("J92").Value = 229.66
$ws.Range("K92").Value = 245.55
$ws.Range("M92").Value = 607.5700000000001
$ws.Range("J93").Value = 205.98
$ws.Range("K93").Value = 242.36
$ws.Range("M93").Value = 640.34
$ws.Range("J94").Value = 229.66
$ws.Range("K94").Value = 245.55
$ws.Range("M94").Value = 607.5700000000001
$ws.Range("J95").Value = 205.9
$ws.Range("K95").Value = 242.33
$ws.Range("M95").Value = 640.3200000000001
$ws.Range("J96").Value = 205.78
$ws.Range("K96").Value = 242.31
$ws.Range("M96").Value = 640.24
$ws.Range("J97").Value = 205.75
$ws.Range("K97").Value = 242.28
$ws.Range("M97").Value = 640.22
